$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue $ws "D2" "42.643.36"
$ws.Range("E2").Value = "  -1.68%  "

Set-TextValue $ws "D3" "2.239.03"
$ws.Range("E3").Value = "  -1.85%  "

$ws.Range("E4").Value = "  -0.01%  "

Set-TextValue $ws "D5" "114.96"
$ws.Range("E5").Value = "  +2.06%  "

Set-TextValue $ws "D6" "285.74"
$ws.Range("E6").Value = "  +7.59%  "

Set-TextValue $ws "D7" "0.626"
$ws.Range("E7").Value = "  -3.79%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("E9").Value = "  +0.11%  "

Set-TextValue $ws "D10" "46.66"
$ws.Range("E10").Value = "  -0.49%  "

Set-TextValue $ws "D11" "0.0930"
$ws.Range("E11").Value = "  -0.56%  "

Set-TextValue $ws "D12" "9.14"
$ws.Range("E12").Value = "  -1.98%  "

$ws.Range("E13").Value = "  -3.14%  "

$ws.Range("E14").Value = "  +0.47%  "

Set-TextValue $ws "D15" "0.883"
$ws.Range("E15").Value = "  +2.22%  "

Set-TextValue $ws "D16" "2.577.79"
$ws.Range("E16").Value = "  -1.75%  "

Set-TextValue $ws "D17" "2.242.71"
$ws.Range("E17").Value = "  -1.55%  "

Set-TextValue $ws "D18" "42.761.98"
$ws.Range("E18").Value = "  -1.37%  "

$ws.Range("E19").Value = "  -1.06%  "

Set-TextValue $ws "D20" "6.93"
$ws.Range("E20").Value = "  +2.61%  "

Set-TextValue $ws "D21" "73.47"
$ws.Range("E21").Value = "  +1.70%  "

$ws.Range("E22").Value = "  +9.66%  "

$ws.Range("E23").Value = "  -3.11%  "

Set-TextValue $ws "D24" "232.27"
$ws.Range("E24").Value = "  -1.15%  "

Set-TextValue $ws "D25" "9.21"
$ws.Range("E25").Value = "  -2.25%  "

Set-TextValue $ws "D26" "12.11"
$ws.Range("E26").Value = "  +5.87%  "

$ws.Range("E27").Value = "  -1.54%  "

Set-TextValue $ws "D28" "3.94"
$ws.Range("E28").Value = "  -0.85%  "

Set-TextValue $ws "D29" "40.09"

$ws.Range("E30").Value = "  -1.47%  "

$ws.Range("E31").Value = "  -0.43%  "

Set-TextValue $ws "D32" "175.64"
$ws.Range("E32").Value = "  +1.28%  "

Set-TextValue $ws "D33" "21.20"
$ws.Range("E33").Value = "  -2.04%  "

Set-TextValue $ws "D34" "0.0905"
$ws.Range("E34").Value = "  +0.88%  "

$ws.Range("E35").Value = "  +18.42%  "

Set-TextValue $ws "D36" "5.60"
$ws.Range("E36").Value = "  -0.83%  "

$ws.Range("E37").Value = "  -3.14%  "

$ws.Range("E38").Value = "  -2.09%  "

$ws.Range("E39").Value = "  -1.45%  "

$ws.Range("E40").Value = "  +1.14%  "

$ws.Range("E41").Value = "  +1.24%  "

Set-TextValue $ws "D42" "72.89"
$ws.Range("E42").Value = "  -2.39%  "

Set-TextValue $ws "D43" "13.53"
$ws.Range("E43").Value = "  -5.32%  "

$ws.Range("E44").Value = "  -1.85%  "

$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("E46").Value = "  -2.26%  "

Set-TextValue $ws "D47" "5.61"
$ws.Range("E47").Value = "  -7.74%  "

$ws.Range("E48").Value = "  +2.13%  "

$ws.Range("E49").Value = "  -0.17%  "

Set-TextValue $ws "D50" "0.650"
$ws.Range("E50").Value = "  +5.86%  "

Set-TextValue $ws "D51" "0.473"
$ws.Range("E51").Value = "  +7.92%  "
